# Adds "sawblade holly berry" cider & liquor recipe rows to the
# "Infused Beverages" sheet, mirroring the existing cider/liquor rows
# (e.g. glass_peach_cider / glass_peach_liquor), and renames the shared
# "toughness" player-choice effect used by glass_peach recipes to
# "projectile_rebound".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Infused Beverages")

# Insert a new row for "sawblade_holly_berry_cider" right after the
# glass_peach_cider row (52), shifting everything below down by one.
$ws.Rows("53:53").Insert()

# Insert a new row for "sawblade_holly_berry_liquor" right after the
# glass_peach_liquor row (60, post first insert), shifting everything
# below down by one more.
$ws.Rows("61:61").Insert()

# Fill the new "sawblade_holly_berry_liquor" row (61) first, then the
# new "sawblade_holly_berry_cider" row (53), so new shared-string
# entries land in the same order as the source edit.
$ws.Range("A61").Value2 = "sawblade_holly_berry_liquor"
$ws.Range("C61").Value2 = "speed"
$ws.Range("F61").Value2 = "slowness"
$ws.Range("G61").Value2 = "nausea"

$ws.Range("A53").Value2 = "sawblade_holly_berry_cider"
$ws.Range("C53").Value2 = "strength"
$ws.Range("F53").Value2 = "weakness"
$ws.Range("G53").Value2 = "nausea"

$ws.Range("I61").Value2 = "sawblade_holly_berry, sugar"
$ws.Range("I53").Value2 = "sawblade_holly_berry"

# Both the glass_peach recipes (rows 52 & 60) and the new sawblade
# recipes (rows 53 & 61) share the "toughness" player-choice column;
# rename that shared effect to "projectile_rebound" for the
# glass_peach rows, and set the new rows' B column to "toughness"
# (re-using the original shared string so indices line up).
$ws.Range("B52").Value2 = "projectile_rebound"
$ws.Range("B60").Value2 = "projectile_rebound"
$ws.Range("B61").Value2 = "toughness"
$ws.Range("B53").Value2 = "toughness"

# Restore the active selection to the newly added liquor row's
# ingredient cell, matching the author's final cursor position.
$ws.Range("I61").Select()
